{"js": "// Replace the 25 division-problem answers in the table with their new values.\n// Each old text is unique within the document, so a direct search+replace is safe.\nconst replacements = [\n  [\"504\u00f73=168, 0\", \"173\u00f77=24, 5\"],\n  [\"551\u00f78=68, 7\", \"328\u00f72=164, 0\"],\n  [\"285\u00f72=142, 1\", \"941\u00f79=104, 5\"],\n  [\"756\u00f79=84, 0\", \"880\u00f76=146, 4\"],\n  [\"794\u00f73=264, 2\", \"643\u00f77=91, 6\"],\n  [\"826\u00f78=103, 2\", \"726\u00f77=103, 5\"],\n  [\"254\u00f72=127, 0\", \"227\u00f75=45, 2\"],\n  [\"398\u00f78=49, 6\", \"107\u00f73=35, 2\"],\n  [\"222\u00f76=37, 0\", \"986\u00f78=123, 2\"],\n  [\"620\u00f72=310, 0\", \"930\u00f75=186, 0\"],\n  [\"499\u00f77=71, 2\", \"824\u00f77=117, 5\"],\n  [\"892\u00f76=148, 4\", \"827\u00f76=137, 5\"],\n  [\"277\u00f78=34, 5\", \"753\u00f74=188, 1\"],\n  [\"944\u00f72=472, 0\", \"497\u00f72=248, 1\"],\n  [\"383\u00f72=191, 1\", \"220\u00f78=27, 4\"],\n  [\"519\u00f79=57, 6\", \"269\u00f74=67, 1\"],\n  [\"227\u00f79=25, 2\", \"681\u00f75=136, 1\"],\n  [\"379\u00f78=47, 3\", \"364\u00f72=182, 0\"],\n  [\"468\u00f73=156, 0\", \"804\u00f73=268, 0\"],\n  [\"120\u00f74=30, 0\", \"453\u00f74=113, 1\"],\n  [\"530\u00f75=106, 0\", \"175\u00f75=35, 0\"],\n  [\"749\u00f74=187, 1\", \"437\u00f79=48, 5\"],\n  [\"612\u00f72=306, 0\", \"448\u00f73=149, 1\"],\n  [\"738\u00f73=246, 0\", \"425\u00f73=141, 2\"],\n  [\"979\u00f73=326, 1\", \"509\u00f79=56, 5\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Replace the 25 division-problem answers in the table with their new values.\n# Each old text is unique within the document, so Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"504\u00f73=168, 0\", \"173\u00f77=24, 5\"),\n    @(\"551\u00f78=68, 7\", \"328\u00f72=164, 0\"),\n    @(\"285\u00f72=142, 1\", \"941\u00f79=104, 5\"),\n    @(\"756\u00f79=84, 0\", \"880\u00f76=146, 4\"),\n    @(\"794\u00f73=264, 2\", \"643\u00f77=91, 6\"),\n    @(\"826\u00f78=103, 2\", \"726\u00f77=103, 5\"),\n    @(\"254\u00f72=127, 0\", \"227\u00f75=45, 2\"),\n    @(\"398\u00f78=49, 6\", \"107\u00f73=35, 2\"),\n    @(\"222\u00f76=37, 0\", \"986\u00f78=123, 2\"),\n    @(\"620\u00f72=310, 0\", \"930\u00f75=186, 0\"),\n    @(\"499\u00f77=71, 2\", \"824\u00f77=117, 5\"),\n    @(\"892\u00f76=148, 4\", \"827\u00f76=137, 5\"),\n    @(\"277\u00f78=34, 5\", \"753\u00f74=188, 1\"),\n    @(\"944\u00f72=472, 0\", \"497\u00f72=248, 1\"),\n    @(\"383\u00f72=191, 1\", \"220\u00f78=27, 4\"),\n    @(\"519\u00f79=57, 6\", \"269\u00f74=67, 1\"),\n    @(\"227\u00f79=25, 2\", \"681\u00f75=136, 1\"),\n    @(\"379\u00f78=47, 3\", \"364\u00f72=182, 0\"),\n    @(\"468\u00f73=156, 0\", \"804\u00f73=268, 0\"),\n    @(\"120\u00f74=30, 0\", \"453\u00f74=113, 1\"),\n    @(\"530\u00f75=106, 0\", \"175\u00f75=35, 0\"),\n    @(\"749\u00f74=187, 1\", \"437\u00f79=48, 5\"),\n    @(\"612\u00f72=306, 0\", \"448\u00f73=149, 1\"),\n    @(\"738\u00f73=246, 0\", \"425\u00f73=141, 2\"),\n    @(\"979\u00f73=326, 1\", \"509\u00f79=56, 5\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $replaced = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $replaced) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\nWrite-Output \"done\""}
